# Plantilla Lista de Tareas de la Entrega 3 - "Termino del CU Consultar informacion maestro"
# Updates the status of several tasks in the "Casos de Uso" sheet to "Hecho"
# and registers consumed hours for days 16 and 17 on the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Task status updates (column F) ---------------------------------------
# Row 23 ("Diagrama entidad-relación"): already "Hecho", just re-center it
# the same way the other already-"Hecho" rows are formatted.
$ws.Range("F23").HorizontalAlignment = -4108

# Row 24 ("Modelo de dominio"): "Por iniciar" -> "Hecho"
$ws.Range("F24").Value = "Hecho"
$ws.Range("F24").HorizontalAlignment = -4108

# Row 26 ("Diagrama entidad-relación"): "En proceso" -> "Hecho"
$ws.Range("F26").Value = "Hecho"

# Row 27 ("Modelo de dominio"): "Por iniciar" -> "Hecho"
$ws.Range("F27").Value = "Hecho"

# --- Consumed hours registered for Dia 16 (AZ) / Dia 17 (BB) --------------
$ws.Range("AZ24").Value = 0.5
$ws.Range("BB24").Value = 1

$ws.Range("AZ26").Value = 6

$ws.Range("AZ27").Value = 0.5
$ws.Range("BB27").Value = 1

# --- Leave the view/selection on the last touched cell ---------------------
$ws.Range("BH24").Select() | Out-Null
